{"js": "// Add a new \"May\" list item after the existing \"April\" list item,\n// matching the same ListParagraph style / numbering as the other\n// month entries.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\n// Find the \"April\" paragraph (last month entry in the list).\nlet aprilParagraph = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text.trim() === \"April\") {\n    aprilParagraph = paragraphs.items[i];\n  }\n}\n\nif (!aprilParagraph) {\n  // Fallback: last paragraph in the body.\n  aprilParagraph = paragraphs.items[paragraphs.items.length - 1];\n}\n\n// Inserting a paragraph after it copies the paragraph's formatting\n// (style + numPr), the same as pressing Enter at the end of \"April\"\n// in Word.\nconst mayParagraph = aprilParagraph.insertParagraph(\"May\", Word.InsertLocation.after);\n\nawait context.sync();\n", "ps1": "# Add a new \"May\" list item after the existing \"April\" list item,\n# matching the same ListParagraph style / numbering as the other\n# month entries in the list.\n$d = $word.ActiveDocument\n\n# Locate the \"April\" paragraph by its index (falls back to the\n# document's last paragraph if, for some reason, no exact match is\n# found).\n$aprilIndex = -1\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    if ($d.Paragraphs.Item($i).Range.Text.Trim() -eq \"April\") {\n        $aprilIndex = $i\n    }\n}\nif ($aprilIndex -eq -1) {\n    $aprilIndex = $d.Paragraphs.Count\n}\n\n$april = $d.Paragraphs.Item($aprilIndex)\n\n# Insert a new paragraph mark right after \"April\" \u2014 like pressing Enter\n# at the end of that line in Word, this carries over the paragraph's\n# formatting (ListParagraph style + list numbering).\n$april.Range.InsertParagraphAfter()\n\n# Structural edits can invalidate earlier handles, so re-fetch the new\n# paragraph from the (now one-longer) Paragraphs collection by index\n# rather than reusing $april.\n$may = $d.Paragraphs.Item($aprilIndex + 1)\n$may.Range.Text = \"May\"\n"}
